$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# TC4 table ("RF008") is missing steps 6/7/8 (only 5 steps present, rows
# 34-38), while the table that follows it (rows 41-53) is pushed straight up
# against it. Bring TC4 up to the same 8-step shape as the other tables by:
#   1) inserting 4 blank rows to make room (pushes the next table, and the
#      merged header cells inside it, down from rows 41-53 to rows 45-57);
#   2) populating the freed rows 39 and 41/42 with steps 6/7/8;
#   3) fixing up the step-5 row (38), which used to hold what is now step 8's
#      text.
# ---------------------------------------------------------------------------

# Make room: push everything from row 41 down by 4 rows.
$ws.Range("A41:F44").Insert()

# --- New row 42 = step 8 ("Avaliador ... 'Salvar'" / "SYSTEM ... salvar...")
#     Grab format + values from row 38, which still holds that exact text
#     (it's about to be overwritten with the real step-5 text below).
$ws.Range("A38:F38").Copy()
$ws.Range("A42:F42").PasteSpecial(-4122)
$ws.Range("A38:F38").Copy()
$ws.Range("A42:F42").PasteSpecial(-4163)
$ws.Range("A42").Value = 8

# --- Row 38 becomes the real step 5 ("Excluir" / "limpa os campos...")
$ws.Range("B38").Value = "Avaliador de Pessoas clica na opcao 'Excluir' para cada Perfil de Competencias avaliado"
$ws.Range("D38").Value = "SYSTEM limpa os campos apresentados 'Nivel da Competencia' e 'Apontamentos' apresentados na tela para cada Perfil de Competencias avaliado"

# --- New row 39 = step 6 ("verifica que os campos ... limpos"), 3 cells only
$ws.Range("A15:C15").Copy()
$ws.Range("A39:C39").PasteSpecial(-4122)
$ws.Range("A15:C15").Copy()
$ws.Range("A39:C39").PasteSpecial(-4163)

# --- New row 41 = step 7 ("Apontamentos" / "preenchidos corretamente")
$ws.Range("A17:F17").Copy()
$ws.Range("A41:F41").PasteSpecial(-4122)
$ws.Range("A17:F17").Copy()
$ws.Range("A41:F41").PasteSpecial(-4163)
